$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 13849.552
$ws.Range("I15").Value = 13849.552
$ws.Range("K15").Value = 41548.656
$ws.Range("M15").Value = -41379.656

$ws.Range("H40").Value = 14977.667
$ws.Range("I40").Value = 19166.666
$ws.Range("J40").Value = 6599.6665
$ws.Range("K40").Value = 19166.666
$ws.Range("L40").Value = 6599.6665
$ws.Range("M40").Value = -18991.666
$ws.Range("N40").Value = -6949.6665

$ws.Range("H63").Value = 43000
$ws.Range("J63").Value = 43000
$ws.Range("L63").Value = 43000
$ws.Range("N63").Value = -44248

$ws.Range("H66").Value = 43000
$ws.Range("J66").Value = 43000
$ws.Range("L66").Value = 129000
$ws.Range("N66").Value = -135240

$ws.Range("H103").Value = 582.0476
$ws.Range("J103").Value = 686.3333
$ws.Range("L103").Value = 2058.9999
$ws.Range("N103").Value = -3230.9999

$ws.Range("H116").Value = 17862750
$ws.Range("I116").Value = 41669492
$ws.Range("J116").Value = 7695.5
$ws.Range("K116").Value = 41669492
$ws.Range("L116").Value = 7695.5
$ws.Range("M116").Value = -41666050
$ws.Range("N116").Value = -14579.5

$ws.Range("H121").Value = 5385.2856
$ws.Range("J121").Value = 5385.2856
$ws.Range("L121").Value = 16155.8568
$ws.Range("N121").Value = -19649.8568

$ws.Range("H137").Value = 3403.262
$ws.Range("I137").Value = 4214.8887
$ws.Range("J137").Value = 2794.5417
$ws.Range("K137").Value = 12644.6661
$ws.Range("L137").Value = 8383.625100000001
$ws.Range("M137").Value = -10094.6661
$ws.Range("N137").Value = -13483.6251

$ws.Range("H138").Value = 3854043.2
$ws.Range("I138").Value = 3571.875
$ws.Range("J138").Value = 10014798
$ws.Range("K138").Value = 10715.625
$ws.Range("L138").Value = 30044394
$ws.Range("M138").Value = -5575.625
$ws.Range("N138").Value = -30054674

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6540.472
$ws.Range("I61").Value = 2617.9092
$ws.Range("J61").Value = 12704.5
$ws.Range("K61").Value = 2617.9092
$ws.Range("L61").Value = 12704.5
$ws.Range("M61").Value = -2405.9092
$ws.Range("N61").Value = -13128.5

$ws.Range("H63").Value = 2349.8333
$ws.Range("J63").Value = 1999.5
$ws.Range("L63").Value = 1999.5
$ws.Range("N63").Value = -3371.5

$ws.Range("H66").Value = 2349.8333
$ws.Range("J66").Value = 1999.5
$ws.Range("L66").Value = 9997.5
$ws.Range("N66").Value = -16861.5

$ws.Range("H136").Value = 6540.472
$ws.Range("I136").Value = 2617.9092
$ws.Range("J136").Value = 12704.5
$ws.Range("K136").Value = 7853.7276
$ws.Range("L136").Value = 38113.5
$ws.Range("M136").Value = -5303.7276
$ws.Range("N136").Value = -43213.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5312.9536
$ws.Range("I134").Value = 2105.7856
$ws.Range("K134").Value = 6317.3568
$ws.Range("M134").Value = -3782.3568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5433.2144
$ws.Range("I16").Value = 3133.125
$ws.Range("J16").Value = 8500
$ws.Range("K16").Value = 3133.125
$ws.Range("L16").Value = 8500
$ws.Range("M16").Value = -2846.125
$ws.Range("N16").Value = -9074

$ws.Range("H113").Value = 5433.2144
$ws.Range("I113").Value = 3133.125
$ws.Range("J113").Value = 8500
$ws.Range("K113").Value = 3133.125
$ws.Range("L113").Value = 8500
$ws.Range("M113").Value = -963.125
$ws.Range("N113").Value = -12840

$ws.Range("H134").Value = 4337.9473
$ws.Range("I134").Value = 1997.5641
$ws.Range("J134").Value = 9408.777
$ws.Range("K134").Value = 5992.692300000001
$ws.Range("L134").Value = 28226.331
$ws.Range("M134").Value = -3457.692300000001
$ws.Range("N134").Value = -33296.331

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2480
$ws.Range("I3").Value = 1630
$ws.Range("K3").Value = 4890
$ws.Range("M3").Value = -4778

$ws.Range("H68").Value = 4294.2856
$ws.Range("I68").Value = 2729.3333
$ws.Range("J68").Value = 4721.091
$ws.Range("K68").Value = 8187.999899999999
$ws.Range("L68").Value = 14163.273
$ws.Range("M68").Value = -7376.999899999999
$ws.Range("N68").Value = -15785.273

$ws.Range("H71").Value = 4294.2856
$ws.Range("I71").Value = 2729.3333
$ws.Range("J71").Value = 4721.091
$ws.Range("K71").Value = 24563.9997
$ws.Range("L71").Value = 42489.819
$ws.Range("M71").Value = -20507.9997
$ws.Range("N71").Value = -50601.819

$ws.Range("H107").Value = 25000362
$ws.Range("I107").Value = 414
$ws.Range("J107").Value = 200000000
$ws.Range("K107").Value = 1242
$ws.Range("L107").Value = 600000000
$ws.Range("M107").Value = 678
$ws.Range("N107").Value = -600003840

$ws.Range("H134").Value = 172264.4
$ws.Range("I134").Value = 172264.4
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 516793.2
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -511723.2
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 3665
$ws.Range("I136").Value = 3665
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10995
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5895
$ws.Range("N136").ClearContents()

$ws.Range("H139").Value = 162752.95
$ws.Range("I139").Value = 233639
$ws.Range("J139").Value = 9166.5
$ws.Range("K139").Value = 700917
$ws.Range("L139").Value = 27499.5
$ws.Range("M139").Value = -695777
$ws.Range("N139").Value = -37779.5

$ws.Range("H140").Value = 335653.1
$ws.Range("J140").Value = 3619.8
$ws.Range("L140").Value = 10859.4
$ws.Range("N140").Value = -21219.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 500
$ws.Range("I33").Value = 500
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 500
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = -248
$ws.Range("N33").Value = -1004

$ws.Range("H58").Value = 72630.5
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 72630.5
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 72630.5
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -73184.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 11520.04
$ws.Range("I136").Value = 6871.857
$ws.Range("K136").Value = 20615.571
$ws.Range("M136").Value = -18065.571

$ws.Range("H140").Value = 56309.668
$ws.Range("J140").Value = 56309.668
$ws.Range("L140").Value = 56309.668
$ws.Range("N140").Value = -66669.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 73333.164
$ws.Range("J46").Value = 73333.164
$ws.Range("L46").Value = 73333.164
$ws.Range("N46").Value = -73795.164

$ws.Range("H68").Value = 43000
$ws.Range("J68").Value = 43000
$ws.Range("L68").Value = 43000
$ws.Range("N68").Value = -44622

$ws.Range("H71").Value = 43000
$ws.Range("J71").Value = 43000
$ws.Range("L71").Value = 129000
$ws.Range("N71").Value = -137112

$ws.Range("H132").Value = 15188317
$ws.Range("I132").Value = 17873444
$ws.Range("K132").Value = 53620332
$ws.Range("M132").Value = -53617802

$ws.Range("H134").Value = 73333.164
$ws.Range("J134").Value = 73333.164
$ws.Range("L134").Value = 219999.492
$ws.Range("N134").Value = -225069.492
